$d = $word.ActiveDocument

# Locate the "Witness B" reading-text paragraph: "Text: when april his showers ..."
$target = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "Text: when april his showers*") {
        $target = $p
    }
}

if ($target -ne $null) {
    $pr = $target.Range
    $words = $pr.Words

    # Words: 1="Text" 2=": " 3="when " 4="april " ...
    $colonSpaceWord = $words.Item(2)   # ": "
    $whenWord = $words.Item(3)         # "when "

    # The space right after "Text:" (currently part of the bold "Text: " run).
    $spaceAfterText = $d.Range($colonSpaceWord.End - 1, $colonSpaceWord.End)

    # The space right after "when" (currently part of the non-bold " " run before "april").
    $spaceAfterWhen = $d.Range($whenWord.End - 1, $whenWord.End)

    # Force Word to split the run containing "Text: " into "Text:" + " ",
    # even though both halves stay bold (toggle off/on dirties the run boundary).
    $spaceAfterText.Bold = 0
    $spaceAfterText.Bold = 1

    # Make the space between "when" and "april" bold as well (was not bold before).
    $spaceAfterWhen.Bold = 1
}
